$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budget summaries")

# New row 38: "Total" check row under the Simulated Historical block (rows 27-36).
# Column C gets the numeric style but no value (matches the other section's blank
# "Total" column), columns D:O get the check formula
#   =<col>28-<col>36-<col>32-<col>30-<col>29
# filled across from D to O.
$ws.Range("C38").NumberFormat = "0.00"
$ws.Range("D38:O38").FormulaR1C1 = "=R[-10]C-R[-2]C-R[-6]C-R[-8]C-R[-9]C"

# New row 39: blank spacer row, but D:O still carry the numeric style.
$ws.Range("D39:O39").NumberFormat = "0.00"

# Update the view's selection to K26 (matches the new sheetView selection).
$ws.Activate() | Out-Null
$ws.Range("K26").Select() | Out-Null
